# Generate Report for handback
#
# The nightly handback run has completed for both target locales: the two
# localized content files (3267b182-... and f0b6f808-...) for zh-cn and
# de-de are now back "in sync with en-US". Refresh the status report:
#   - flip each row's Status from "Ready for handoff" to
#     "Handed back: in sync with en-US" (Overview + both locale sheets)
#   - populate "Latest Target File" / "Latest Handback File" for the two
#     localized rows on each locale sheet, with hyperlinks matching the
#     existing Source/Handoff-file links
#   - stamp "Latest Handback DateTime" with the handback timestamp

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: just the status text (columns B and C) -----------
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("B2").Value = $newStatus
$ov.Range("C2").Value = $newStatus
$ov.Range("B3").Value = $newStatus
$ov.Range("C3").Value = $newStatus

# --- zh-cn sheet --------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("B2").Value = $newStatus
$ws.Range("E2").Value = "3267b182-833f-4883-9d8e-c186289f95b4.md"
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/466d159288e7203b14f1516bc9f7ca084aff9ae3/e2e/3267b182-833f-4883-9d8e-c186289f95b4.md", [Type]::Missing, [Type]::Missing, "3267b182-833f-4883-9d8e-c186289f95b4.md") | Out-Null
$ws.Range("F2").Value = "3267b182-833f-4883-9d8e-c186289f95b4.b2efe56a1def11e37d7c1a8ad00ecf2a587e58c8.zh-cn.xlf"
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3d36dd319ef77161fdf955318a0b5ad4b813a731/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/3267b182-833f-4883-9d8e-c186289f95b4.b2efe56a1def11e37d7c1a8ad00ecf2a587e58c8.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "3267b182-833f-4883-9d8e-c186289f95b4.b2efe56a1def11e37d7c1a8ad00ecf2a587e58c8.zh-cn.xlf") | Out-Null
$ws.Range("G2").Value = "2016-01-18 03:06:59"

$ws.Range("B3").Value = $newStatus
$ws.Range("E3").Value = "f0b6f808-0bb2-45d4-96a0-a9774781cba5.md"
$ws.Hyperlinks.Add($ws.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/466d159288e7203b14f1516bc9f7ca084aff9ae3/e2e/f0b6f808-0bb2-45d4-96a0-a9774781cba5.md", [Type]::Missing, [Type]::Missing, "f0b6f808-0bb2-45d4-96a0-a9774781cba5.md") | Out-Null
$ws.Range("F3").Value = "f0b6f808-0bb2-45d4-96a0-a9774781cba5.e0ef89a85990f1a1915e91be30101887265ab26b.zh-cn.xlf"
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3d36dd319ef77161fdf955318a0b5ad4b813a731/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/f0b6f808-0bb2-45d4-96a0-a9774781cba5.e0ef89a85990f1a1915e91be30101887265ab26b.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "f0b6f808-0bb2-45d4-96a0-a9774781cba5.e0ef89a85990f1a1915e91be30101887265ab26b.zh-cn.xlf") | Out-Null
$ws.Range("G3").Value = "2016-01-18 03:06:59"

# --- de-de sheet ---------------------------------------------------------
$ws2 = $wb.Worksheets.Item("de-de")

$ws2.Range("B2").Value = $newStatus
$ws2.Range("E2").Value = "3267b182-833f-4883-9d8e-c186289f95b4.md"
$ws2.Hyperlinks.Add($ws2.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/466d159288e7203b14f1516bc9f7ca084aff9ae3/e2e/3267b182-833f-4883-9d8e-c186289f95b4.md", [Type]::Missing, [Type]::Missing, "3267b182-833f-4883-9d8e-c186289f95b4.md") | Out-Null
$ws2.Range("F2").Value = "3267b182-833f-4883-9d8e-c186289f95b4.b2efe56a1def11e37d7c1a8ad00ecf2a587e58c8.de-de.xlf"
$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0de9d5950e84e42210ba83f2ec374c787086c7b9/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/3267b182-833f-4883-9d8e-c186289f95b4.b2efe56a1def11e37d7c1a8ad00ecf2a587e58c8.de-de.xlf", [Type]::Missing, [Type]::Missing, "3267b182-833f-4883-9d8e-c186289f95b4.b2efe56a1def11e37d7c1a8ad00ecf2a587e58c8.de-de.xlf") | Out-Null
$ws2.Range("G2").Value = "2016-01-18 03:07:21"

$ws2.Range("B3").Value = $newStatus
$ws2.Range("E3").Value = "f0b6f808-0bb2-45d4-96a0-a9774781cba5.md"
$ws2.Hyperlinks.Add($ws2.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/466d159288e7203b14f1516bc9f7ca084aff9ae3/e2e/f0b6f808-0bb2-45d4-96a0-a9774781cba5.md", [Type]::Missing, [Type]::Missing, "f0b6f808-0bb2-45d4-96a0-a9774781cba5.md") | Out-Null
$ws2.Range("F3").Value = "f0b6f808-0bb2-45d4-96a0-a9774781cba5.e0ef89a85990f1a1915e91be30101887265ab26b.de-de.xlf"
$ws2.Hyperlinks.Add($ws2.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0de9d5950e84e42210ba83f2ec374c787086c7b9/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/f0b6f808-0bb2-45d4-96a0-a9774781cba5.e0ef89a85990f1a1915e91be30101887265ab26b.de-de.xlf", [Type]::Missing, [Type]::Missing, "f0b6f808-0bb2-45d4-96a0-a9774781cba5.e0ef89a85990f1a1915e91be30101887265ab26b.de-de.xlf") | Out-Null
$ws2.Range("G3").Value = "2016-01-18 03:07:21"

Write-Output "handback report updated"
